$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E6").Value = "AppleExerciseTime"
$ws.Range("F6").Value = 60
$ws.Range("F7").Select()
